$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("M2").Value = 41.474781
$ws.Range("N2").Value = 82.949562
$ws.Range("O2").Value = 0.6394851352970483
$ws.Range("P2").Value = 0.5499336598827257
$ws.Range("Q2").Value = 1.160478197307
$ws.Range("R2").Value = 6.962869183842
$ws.Range("S2").Value = 0.6394851352970483
$ws.Range("T2").Value = 0.5499336598827257

# Row 3
$ws.Range("M3").Value = 9.567994999999998
$ws.Range("N3").Value = 28.703985
$ws.Range("O3").Value = 0.1475255668522151
$ws.Range("P3").Value = 0.1902998297238611
$ws.Range("Q3").Value = 0.2677156894316666
$ws.Range("R3").Value = 2.409441204885
$ws.Range("S3").Value = 0.1475255668522151
$ws.Range("T3").Value = 0.1902998297238611

# Row 4
$ws.Range("M4").Value = 5.805205666666667
$ws.Range("N4").Value = 17.415617
$ws.Range("O4").Value = 0.0895084348046473
$ws.Range("P4").Value = 0.1154609351153152
$ws.Range("Q4").Value = 0.1624315896218889
$ws.Range("R4").Value = 1.461884306597
$ws.Range("S4").Value = 0.0895084348046473
$ws.Range("T4").Value = 0.1154609351153152

# Row 5
$ws.Range("M5").Value = 2.2592025
$ws.Range("N5").Value = 4.518405
$ws.Range("O5").Value = 0.03483385280264482
$ws.Range("P5").Value = 0.02995583024877705
$ws.Range("Q5").Value = 0.06321323901750001
$ws.Range("R5").Value = 0.379279434105
$ws.Range("S5").Value = 0.03483385280264482
$ws.Range("T5").Value = 0.02995583024877705

# Row 6
$ws.Range("M6").Value = 1.733200666666667
$ws.Range("N6").Value = 5.199602000000001
$ws.Range("O6").Value = 0.02672361459413777
$ws.Range("P6").Value = 0.03447198621487044
$ws.Range("Q6").Value = 0.04849553238688889
$ws.Range("R6").Value = 0.436459791482
$ws.Range("S6").Value = 0.02672361459413777
$ws.Range("T6").Value = 0.03447198621487044

# Row 7
$ws.Range("M7").Value = 4.016136
$ws.Range("N7").Value = 12.048408
$ws.Range("O7").Value = 0.06192339564930666
$ws.Range("P7").Value = 0.07987775881445054
$ws.Range("Q7").Value = 0.112372823992
$ws.Range("R7").Value = 1.011355415928
$ws.Range("S7").Value = 0.06192339564930666
$ws.Range("T7").Value = 0.07987775881445054
